$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '42.637.72'
$ws.Range("E2").Value = '  +0.67%  '
Set-TextValue $ws.Range("D3") '2.300.46'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.08%  '
Set-TextValue $ws.Range("D5") '316.39'
$ws.Range("E5").Value = '  +0.18%  '
Set-TextValue $ws.Range("D6") '103.66'
$ws.Range("E6").Value = '  -0.73%  '
Set-TextValue $ws.Range("D7") '0.628'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").Value = '  +0.09%  '
Set-TextValue $ws.Range("D9") '0.605'
$ws.Range("E9").Value = '  -0.44%  '
Set-TextValue $ws.Range("D10") '39.79'
$ws.Range("E10").Value = '  +0.25%  '
Set-TextValue $ws.Range("D11") '0.0904'
$ws.Range("E11").Value = '  -0.43%  '
Set-TextValue $ws.Range("D12") '8.51'
$ws.Range("E12").Value = '  +2.39%  '
$ws.Range("E13").Value = '  +0.69%  '
Set-TextValue $ws.Range("D14") '0.992'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("E15").Value = '  +0.21%  '
Set-TextValue $ws.Range("D16") '2.652.09'
$ws.Range("E16").Value = '  +0.18%  '
Set-TextValue $ws.Range("D17") '2.303.17'
$ws.Range("E17").Value = '  +0.22%  '
Set-TextValue $ws.Range("D18") '42.573.45'
$ws.Range("E18").Value = '  +0.69%  '
Set-TextValue $ws.Range("D19") '7.63'
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("E20").Value = '  +0.48%  '
Set-TextValue $ws.Range("D21") '13.77'
$ws.Range("E21").Value = '  +31.77%  '
Set-TextValue $ws.Range("D22") '74.01'
$ws.Range("E22").Value = '  +0.75%  '
Set-TextValue $ws.Range("D23") '3.54'
$ws.Range("E23").Value = '  -1.80%  '
Set-TextValue $ws.Range("D24") '267.45'
$ws.Range("E24").Value = '  -3.92%  '
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  +0.34%  '
Set-TextValue $ws.Range("D28") '2.27'
$ws.Range("E28").Value = '  -3.92%  '
Set-TextValue $ws.Range("D29") '22.60'
$ws.Range("E29").Value = '  -1.18%  '
Set-TextValue $ws.Range("D30") '6.64'
$ws.Range("E30").Value = '  +13.85%  '
Set-TextValue $ws.Range("D31") '37.61'
$ws.Range("E31").Value = '  +3.61%  '
Set-TextValue $ws.Range("D32") '165.40'
$ws.Range("E32").Value = '  +0.71%  '
Set-TextValue $ws.Range("D33") '0.0883'
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("E34").Value = '  -5.29%  '
Set-TextValue $ws.Range("D35") '0.131'
$ws.Range("E35").Value = '  -3.21%  '
Set-TextValue $ws.Range("D36") '0.113'
$ws.Range("E36").Value = '  -0.14%  '
Set-TextValue $ws.Range("D37") '4.59'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("E38").Value = '  +1.49%  '
Set-TextValue $ws.Range("D39") '3.72'
$ws.Range("E39").Value = '  -1.36%  '
Set-TextValue $ws.Range("D40") '2.74'
$ws.Range("E40").Value = '  -2.27%  '
Set-TextValue $ws.Range("D41") '1.61'
$ws.Range("E41").Value = '  +10.88%  '
Set-TextValue $ws.Range("D42") '98.03'
$ws.Range("E42").Value = '  -1.48%  '
Set-TextValue $ws.Range("D43") '70.06'
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("E45").Value = '  -0.07%  '
Set-TextValue $ws.Range("D46") '12.38'
$ws.Range("E46").Value = '  +2.93%  '
Set-TextValue $ws.Range("D47") '116.58'
$ws.Range("E47").Value = '  +3.94%  '
Set-TextValue $ws.Range("D48") '80.74'
$ws.Range("E48").Value = '  +3.91%  '
Set-TextValue $ws.Range("D49") '1.638.44'
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D50") '8.92'
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D51") '5.29'
$ws.Range("E51").Value = '  -0.09%  '
